$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 2369050
$ws.Cells.Item(4, 3).Value = 12393
$ws.Cells.Item(4, 4).Value = 984277
$ws.Cells.Item(4, 5).Value = 1262375
$ws.Cells.Item(4, 7).Value = 151
$ws.Cells.Item(4, 8).Value = 122398

# Row 7
$ws.Cells.Item(7, 2).Value = 440183
$ws.Cells.Item(7, 3).Value = 13273
$ws.Cells.Item(7, 4).Value = 248129
$ws.Cells.Item(7, 5).Value = 178039
$ws.Cells.Item(7, 7).Value = 312
$ws.Cells.Item(7, 8).Value = 14015

# Row 14
$ws.Cells.Item(14, 2).Value = 191912
$ws.Cells.Item(14, 3).Value = 337
$ws.Cells.Item(14, 5).Value = 7648

# Row 33
$ws.Cells.Item(33, 2).Value = 45303
$ws.Cells.Item(33, 3).Value = 378
$ws.Cells.Item(33, 4).Value = 33046
$ws.Cells.Item(33, 5).Value = 11954
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 303

# Row 50
$ws.Cells.Item(50, 5).Value = 5280
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 65

# Row 76
$ws.Cells.Item(76, 2).Value = 6461
$ws.Cells.Item(76, 3).Value = 146
$ws.Cells.Item(76, 5).Value = 1992

# Row 79
$ws.Cells.Item(79, 2).Value = 5513
$ws.Cells.Item(79, 3).Value = 56
$ws.Cells.Item(79, 4).Value = 4039
$ws.Cells.Item(79, 5).Value = 1422

# Row 103
$ws.Cells.Item(103, 1).Value = 'Maldivas'
$ws.Cells.Item(103, 2).Value = 2217
$ws.Cells.Item(103, 3).Value = 14
$ws.Cells.Item(103, 4).Value = 1813
$ws.Cells.Item(103, 5).Value = 396
$ws.Cells.Item(103, 8).Value = 8

# Row 104
$ws.Cells.Item(104, 1).Value = 'Costa Rica'
$ws.Cells.Item(104, 2).Value = 2213
$ws.Cells.Item(104, 4).Value = 1032
$ws.Cells.Item(104, 5).Value = 1169
$ws.Cells.Item(104, 8).Value = 12

# Row 129
$ws.Cells.Item(129, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(129, 2).Value = 989
$ws.Cells.Item(129, 3).Value = 156
$ws.Cells.Item(129, 4).Value = 442
$ws.Cells.Item(129, 5).Value = 544
$ws.Cells.Item(129, 8).Value = 3

# Row 130
$ws.Cells.Item(130, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(130, 2).Value = 988
$ws.Cells.Item(130, 3).Value = 2
$ws.Cells.Item(130, 4).Value = 824
$ws.Cells.Item(130, 5).Value = 145
$ws.Cells.Item(130, 8).Value = 19

# Row 131
$ws.Cells.Item(131, 2).Value = 967
$ws.Cells.Item(131, 3).Value = 26
$ws.Cells.Item(131, 4).Value = 350
$ws.Cells.Item(131, 5).Value = 360
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 257

# Row 132
$ws.Cells.Item(132, 1).Value = 'Cabo Verde'
$ws.Cells.Item(132, 2).Value = 944
$ws.Cells.Item(132, 3).Value = 54
$ws.Cells.Item(132, 4).Value = 419
$ws.Cells.Item(132, 5).Value = 517
$ws.Cells.Item(132, 8).Value = 8

# Row 133
$ws.Cells.Item(133, 1).Value = 'Georgia'
$ws.Cells.Item(133, 2).Value = 908
$ws.Cells.Item(133, 3).Value = 2
$ws.Cells.Item(133, 4).Value = 761
$ws.Cells.Item(133, 5).Value = 133
$ws.Cells.Item(133, 8).Value = 14

# Row 134
$ws.Cells.Item(134, 1).Value = 'Burkina Faso'
$ws.Cells.Item(134, 2).Value = 903
$ws.Cells.Item(134, 4).Value = 814
$ws.Cells.Item(134, 5).Value = 36
$ws.Cells.Item(134, 8).Value = 53

# Row 141
$ws.Cells.Item(141, 2).Value = 737
$ws.Cells.Item(141, 3).Value = 4
$ws.Cells.Item(141, 5).Value = 551

# Row 144
$ws.Cells.Item(144, 2).Value = 702
$ws.Cells.Item(144, 3).Value = 4
$ws.Cells.Item(144, 5).Value = 487

# Row 148
$ws.Cells.Item(148, 1).Value = 'Liberia'
$ws.Cells.Item(148, 2).Value = 650
$ws.Cells.Item(148, 3).Value = 24
$ws.Cells.Item(148, 4).Value = 260
$ws.Cells.Item(148, 5).Value = 356
$ws.Cells.Item(148, 8).Value = 34

# Row 149
$ws.Cells.Item(149, 1).Value = 'Suazilandia'
$ws.Cells.Item(149, 2).Value = 643
$ws.Cells.Item(149, 3).Value = 8
$ws.Cells.Item(149, 4).Value = 291
$ws.Cells.Item(149, 5).Value = 346
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 6

# Row 156
$ws.Cells.Item(156, 2).Value = 367
$ws.Cells.Item(156, 3).Value = 5
$ws.Cells.Item(156, 5).Value = 43

# Row 168
$ws.Cells.Item(168, 1).Value = 'Angola'
$ws.Cells.Item(168, 2).Value = 186
$ws.Cells.Item(168, 3).Value = 3
$ws.Cells.Item(168, 4).Value = 77
$ws.Cells.Item(168, 5).Value = 99
$ws.Cells.Item(168, 7).Value = 1
$ws.Cells.Item(168, 8).Value = 10

# Row 169
$ws.Cells.Item(169, 1).Value = 'Guyana'
$ws.Cells.Item(169, 2).Value = 184
$ws.Cells.Item(169, 4).Value = 103
$ws.Cells.Item(169, 5).Value = 69
$ws.Cells.Item(169, 8).Value = 12

# Row 202
$ws.Cells.Item(202, 1).Value = 'Dominica'

# Row 203
$ws.Cells.Item(203, 1).Value = 'Fiyi'

# Update timestamp text in A1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 22 de Junio de 2020 a las 20:44'
